# "config to run all weather years"
# Applies the edits to the "Coupling Parameters" sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# --- Row 15: npv calculated with annuity instead of restpayment -----------
# Update the FALSE-branch text of the formula (adds a "don't use this" note)
# before flipping the switch, so the cached text lines up either way.
$ws.Range("C15").Formula = '=IF(B15=TRUE,"the npv is calculated with the annuity","the npv is calculated with the restpayment _ >don' + [char]0x2019 + 't use this")'
$ws.Range("B15").Value = $true
$ws.Range("B15").Interior.Color = 65535

# --- Row 23: dummy_capacity_to_test 0 -> 100, highlighted ------------------
$ws.Range("B23").Value = 100
$ws.Range("B23").Interior.Color = 65535

# --- Row 28 & 29: stop fixing demand/profiles to representative year ------
$ws.Range("B28").Value = $false
$ws.Range("B29").Value = $false
$ws.Range("B28:B29").Interior.Color = 65535

# --- Row 30: C30 becomes a formula instead of a fixed explanatory string --
$ws.Range("C30").Formula = '=IF(AND(B28=FALSE,B29=FALSE),"NOTSET","if NOTSET then future year considers look ahead. Otherwise it considers this future year")'

# --- View state: scrolled down a bit, new selected cell -------------------
$ws.Application.Goto($ws.Range("A18"))
$ws.Range("C20").Select()

$wb.Application.Calculate()
